# Refresh the cryptocurrency price (D) / 1h volume change (E) snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.858.75"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.887.19"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.7462"
$ws.Range("E5").Value = "  -4.69%  "
$ws.Range("D6").Value = "'242.74"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.3113"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").Value = "'25.39"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'0.07118"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").Value = "'0.08473"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").Value = "'0.7599"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "1.906.72"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'5.353"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'93.27"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "'6.147"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "29.990.08"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "'243.35"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "'0.000007792"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "2.160.56"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'8.011"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'0.1590"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'9.381"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "'162.49"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'2.026"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'1.513"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'4.469"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'4.101"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "'0.05397"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("D35").Value = "'1.239"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "'0.7440"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'2.712"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "'0.01929"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'2.766"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'0.4449"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'6.071"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "'72.64"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "1.083.08"
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("D45").Value = "'0.8622"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'102.58"
$ws.Range("D48").Value = "'7.657"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "'1.858"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "'3.056"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("D51").Value = "2.055.28"
$ws.Range("E51").Value = "  +1.35%  "

# Drop the quote-prefix style picked up above so cell styling stays untouched.
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
